$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8:T10").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spon2"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8164013333333333
$ws.Range("H2").Value = 2.449204
$ws.Range("I2").Value = 0.05618115571687973
$ws.Range("J2").Value = 0.05618115571687973
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005673666666666667
$ws.Range("N2").Value = 0.017021
$ws.Range("O2").Value = 0.1234828534325781
$ws.Range("P2").Value = 0.1234828534325781
$ws.Range("Q2").Value = 0.004631989031555555
$ws.Range("R2").Value = 0.041687901284
$ws.Range("S2").Value = 0.006937409417060308
$ws.Range("T2").Value = 0.006937409417060309

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spon2"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8164013333333333
$ws.Range("H3").Value = 2.449204
$ws.Range("I3").Value = 0.05618115571687973
$ws.Range("J3").Value = 0.05618115571687973
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04027333333333333
$ws.Range("N3").Value = 0.12082
$ws.Range("O3").Value = 0.8765171465674219
$ws.Range("P3").Value = 0.876517146567422
$ws.Range("Q3").Value = 0.03287920303111111
$ws.Range("R3").Value = 0.29591282728
$ws.Range("S3").Value = 0.04924374629981942
$ws.Range("T3").Value = 0.04924374629981943

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Spon2"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.878362
$ws.Range("H4").Value = 38.635086
$ws.Range("I4").Value = 0.8862323361798529
$ws.Range("J4").Value = 0.8862323361798529
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005673666666666667
$ws.Range("N4").Value = 0.017021
$ws.Range("O4").Value = 0.1234828534325781
$ws.Range("P4").Value = 0.1234828534325781
$ws.Range("Q4").Value = 0.07306753320066668
$ws.Range("R4").Value = 0.6576077988060001
$ws.Range("S4").Value = 0.1094344976757081
$ws.Range("T4").Value = 0.1094344976757081

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Spon2"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.878362
$ws.Range("H5").Value = 38.635086
$ws.Range("I5").Value = 0.8862323361798529
$ws.Range("J5").Value = 0.8862323361798529
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04027333333333333
$ws.Range("N5").Value = 0.12082
$ws.Range("O5").Value = 0.8765171465674219
$ws.Range("P5").Value = 0.876517146567422
$ws.Range("Q5").Value = 0.5186545656133333
$ws.Range("R5").Value = 4.66789109052
$ws.Range("S5").Value = 0.7767978385041447
$ws.Range("T5").Value = 0.7767978385041449

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Spon2"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8368233333333334
$ws.Range("H6").Value = 2.51047
$ws.Range("I6").Value = 0.05758650810326746
$ws.Range("J6").Value = 0.05758650810326746
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005673666666666667
$ws.Range("N6").Value = 0.017021
$ws.Range("O6").Value = 0.1234828534325781
$ws.Range("P6").Value = 0.1234828534325781
$ws.Range("Q6").Value = 0.004747856652222222
$ws.Range("R6").Value = 0.04273070987000001
$ws.Range("S6").Value = 0.007110946339809748
$ws.Range("T6").Value = 0.007110946339809748

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Spon2"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8368233333333334
$ws.Range("H7").Value = 2.51047
$ws.Range("I7").Value = 0.05758650810326746
$ws.Range("J7").Value = 0.05758650810326746
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04027333333333333
$ws.Range("N7").Value = 0.12082
$ws.Range("O7").Value = 0.8765171465674219
$ws.Range("P7").Value = 0.876517146567422
$ws.Range("Q7").Value = 0.03370166504444445
$ws.Range("R7").Value = 0.3033149854
$ws.Range("S7").Value = 0.05047556176345771
$ws.Range("T7").Value = 0.05047556176345772
